$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''65.140.44'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  +3.16%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''2.631.02'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  +2.00%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  +0.08%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''596.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  +1.61%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''155.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  +5.08%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = '''  +0.05%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = '''  +1.12%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = '''  +8.88%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''0.402'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  +5.51%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''5.79'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  +1.06%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = '''  +2.07%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''29.14'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  +6.83%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''0.0000186'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  +22.35%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''3.103.68'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  +2.19%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''65.014.06'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  +3.30%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''2.628.79'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  +1.89%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''12.54'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  +3.60%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D20').Value = '''352.04'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  +2.48%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = '''  +8.46%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = '''  +0.24%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''68.43'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  +2.77%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''9.52'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  +5.06%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = '''  -3.19%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = '''  -0.38%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = '''  +1.52%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = '''  +0.58%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''0.0₃0952'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  +12.72%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = '''  +0.02%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''528.32'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -4.46%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''2.11'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  +5.02%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = '''  +2.53%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = '''  +7.98%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = '''  +6.75%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = '''  +3.89%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''20.30'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  +5.20%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''163.76'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -0.97%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = '''  +5.94%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = '''  -0.01%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = '''  -0.04%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''42.30'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +7.12%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''165.25'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -0.01%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''4.10'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  +4.05%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''0.0617'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  +5.55%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''23.06'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  +2.21%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = '''  +9.14%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = '''  +3.13%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = '''  +3.63%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = '''  +2.11%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''19.42'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  +2.97%  '
$ws.Range('E51').Style = 'Normal'
